$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2014 -> 01/01/2021 (keep as text, like the original cell)
$ws.Range("B8").Value = "'01/01/2021"
$ws.Range("C8").Value = "'01/01/2021"

# Docentes responsáveis: Marco Antonio Carvalho Pereira -> Herlandí de Souza Andrade
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# Método:
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"

# Critério:
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# Norma de recuperação:
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
